$d = $word.ActiveDocument

# --- Fix 1: remove stray trailing "a" from footnote text "...པེ་ཅིན།a" -> "...པེ་ཅིན།"
$targetFootnote = $null
for ($i = 1; $i -le $d.Footnotes.Count; $i++) {
    $fn = $d.Footnotes.Item($i)
    $t = $fn.Range.Text
    if ($t.Length -gt 0 -and $t.EndsWith("a")) {
        $targetFootnote = $fn
        break
    }
}
if ($targetFootnote -ne $null) {
    $full = $targetFootnote.Range.Text
    $trimmed = $full.Substring(0, $full.Length - 1)
    $targetFootnote.Range.Text = $trimmed
}

# --- Fix 2: delete the trailing empty footnote (body text is just "।", i.e. a
# single punctuation mark with no real annotation) and its reference in the
# document body, which resolves the "bug with empty notes".
for ($i = $d.Footnotes.Count; $i -ge 1; $i--) {
    $fn = $d.Footnotes.Item($i)
    $txt = $fn.Range.Text
    $trimmedTxt = $txt.Trim()
    if ($trimmedTxt.Length -le 1) {
        $fn.Delete()
    }
}
